$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("case 1")

$ws.Rows("15:16").Insert()

$ws.Range("A15").Value = "GSG"
$ws.Range("B15").Value = "state"
$ws.Range("C15").Value = ""

$ws.Range("A15:C15").Interior.Pattern = -4142

$ws.Range("B15").Borders.Item(10).LineStyle = 1
$ws.Range("B15").Borders.Item(10).Weight = 4
$ws.Range("B15").Borders.Item(10).ColorIndex = 1

$ws.Range("B16").Value = "goal sig"
$ws.Range("B16").Interior.Color = $ws.Range("B14").Interior.Color
$ws.Range("B16").Borders.Item(10).LineStyle = $ws.Range("B14").Borders.Item(10).LineStyle
$ws.Range("B16").Borders.Item(10).Weight = $ws.Range("B14").Borders.Item(10).Weight
$ws.Range("B16").Borders.Item(10).ColorIndex = $ws.Range("B14").Borders.Item(10).ColorIndex

$shp = $ws.Shapes.Item(1)
$shp.Height = $shp.Height + (14.4 * 2)

$ws.Range("E15").Select()

$wb.Save()
